$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (Project Start Date / Project End Date) before
# the old "Type" column (I), shifting old I:T -> K:V.
$ws.Columns("I:J").Insert()

# Row 1 - new column headers
$ws.Range("I1").Value = "Project Start Date"
$ws.Range("J1").Value = "Project End Date"

# Row 2 - field type hints
$ws.Range("I2").Value = "Date field"
$ws.Range("J2").Value = "Date field"

# Row 3 - help text (style is carried over from the old I3 cell by the insert)
$ws.Range("I3").Value = "If possible please provide dates in format yyyy-mm-dd"
$ws.Range("J3").Value = "If possible please provide dates in format yyyy-mm-dd"

# Row 4 - extra note under the new End Date column
$ws.Range("J4").Value = "Approximate/estimated end dates accepted"

# Update the active selection to reflect where editing left off
$null = $ws.Range("J3").Select()
